$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 249, shifting existing rows 249:352 down to 250:353
$ws.Rows("249:249").Insert()

# Populate the new row 249 with the new data record
$ws.Range("A249").Value = 6
$ws.Range("B249").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 45141
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112029
$ws.Range("G249").Value = "Orégano"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 28
$ws.Range("K249").Value = 20000
$ws.Range("L249").Value = 20000
$ws.Range("M249").Value = 20000
$ws.Range("N249").Value = "$/docena de atados"
$ws.Range("O249").Value = "Región Metropolitana"
$ws.Range("P249").Value = 6667
$ws.Range("Q249").Value = 3
$ws.Range("R249").Value = "Hortaliza"
